$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.101398468017578
$ws.Range("B1").Value = 4.530810356140137
$ws.Range("C1").Value = 5.696710586547852
$ws.Range("D1").Value = 8.268106460571289
$ws.Range("E1").Value = 4.407927989959717
